$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Header label: "LCSC" -> "LCSC#"
# ------------------------------------------------------------------
$ws.Range("L8").Value = "LCSC#"

# ------------------------------------------------------------------
# 2. Summary block (rows 2-6): regenerated KiBot run with KiCad 7.0.8
#    and a smaller BOM (3 groups / 3 THT parts instead of 10).
# ------------------------------------------------------------------
$ws.Range("D6").Value = "7.0.8-7.0.8~ubuntu22.04.1"
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = "3 (0 SMD/ 3 THT)"
$ws.Range("F4").Value = "3 (0 SMD/ 3 THT)"
$ws.Range("F6").Value = 3

# ------------------------------------------------------------------
# 3. BOM table: the first 4 component rows (APA102 LED strip +
#    0.7mm/1.5mm cable connectors) disappear from the regenerated
#    BOM; the Audio Jack (old row 13) and XT30 (old row 14) rows
#    become the new row 9 / row 10.  Stash the "Row" numbering cells
#    first (they already hold the text values "1"/"2" we need), move
#    the surviving data up, then restore those text cells.
# ------------------------------------------------------------------
$ws.Range("A9").Copy($ws.Range("Z1"))
$ws.Range("A10").Copy($ws.Range("Z2"))

$ws.Range("A13:S13").Copy($ws.Range("A9:S9"))
$ws.Range("A14:S14").Copy($ws.Range("A10:S10"))

$ws.Range("Z1").Copy()
$ws.Range("A9").PasteSpecial(-4163)
$ws.Range("Z2").Copy()
$ws.Range("A10").PasteSpecial(-4163)
$ws.Range("Z1:Z2").Clear()

# new row 9 (Audio Jack) no longer uses the 30pt custom row height
# that the old row 9 (and old row 13's blank default) differed on.
$ws.Rows("9").AutoFit()

# drop the now-duplicated old rows 11-14
$ws.Range("A11:A14").EntireRow.Delete()

# ------------------------------------------------------------------
# 4. Column widths: Value column narrower, LCSC# column wider.
# ------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 24.86
$ws.Columns("L").ColumnWidth = 14.8

Write-Output "done"
